$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state (from diff): 12 data rows (rows 2-13), cross join of
# Sending cluster {ECs, FAPs, sCs} x Target cluster {ECs, FAPs, M2, sCs}
# with recomputed NATMI specificity statistics for 3 sending clusters.
$data = @(
    @("ECs", "Myoc", "Fzd7", "ECs", 1, 0.3333333333333333, 0.2144083333333333, 0.643225, 0.008611346839948651, 0.008611346839948651, 3, 1, 2.500933333333333, 7.5028, 0.129395502869887, 0.129395502869887, 0.5362209477777777, 4.82598853, 0.001114269554742168, 0.001114269554742168),
    @("ECs", "Myoc", "Fzd7", "FAPs", 1, 0.3333333333333333, 0.2144083333333333, 0.643225, 0.008611346839948651, 0.008611346839948651, 3, 1, 7.782607, 23.347821, 0.4026634108880828, 0.4026634108880828, 1.668655795858333, 15.017902162725, 0.003467474290914037, 0.003467474290914037),
    @("ECs", "Myoc", "Fzd7", "M2", 1, 0.3333333333333333, 0.2144083333333333, 0.643225, 0.008611346839948651, 0.008611346839948651, 3, 1, 4.565521, 13.696563, 0.2362149673420707, 0.2362149673420707, 0.9788857484083334, 8.809971735675001, 0.002034129012569714, 0.002034129012569714),
    @("ECs", "Myoc", "Fzd7", "sCs", 1, 0.3333333333333333, 0.2144083333333333, 0.643225, 0.008611346839948651, 0.008611346839948651, 3, 1, 4.478761333333333, 13.436284, 0.2317261188999595, 0.2317261188999595, 0.9602837528777777, 8.6425537759, 0.001995473981722732, 0.001995473981722732),
    @("FAPs", "Myoc", "Fzd7", "ECs", 3, 1, 23.971258, 71.91377399999999, 0.9627648963950115, 0.9627648963950115, 3, 1, 2.500933333333333, 7.5028, 0.129395502869887, 0.129395502869887, 59.95051817413332, 539.5546635671999, 0.1245774479145072, 0.1245774479145071),
    @("FAPs", "Myoc", "Fzd7", "FAPs", 3, 1, 23.971258, 71.91377399999999, 0.9627648963950115, 0.9627648963950115, 3, 1, 7.782607, 23.347821, 0.4026634108880828, 0.4026634108880828, 186.558880309606, 1679.029922786454, 0.387670197065727, 0.387670197065727),
    @("FAPs", "Myoc", "Fzd7", "M2", 3, 1, 23.971258, 71.91377399999999, 0.9627648963950115, 0.9627648963950115, 3, 1, 4.565521, 13.696563, 0.2362149673420707, 0.2362149673420707, 109.441281795418, 984.971536158762, 0.2274194785600397, 0.2274194785600397),
    @("FAPs", "Myoc", "Fzd7", "sCs", 3, 1, 23.971258, 71.91377399999999, 0.9627648963950115, 0.9627648963950115, 3, 1, 4.478761333333333, 13.436284, 0.2317261188999595, 0.2317261188999595, 107.3615434417573, 966.2538909758158, 0.2230977728547377, 0.2230977728547377),
    @("sCs", "Myoc", "Fzd7", "ECs", 3, 1, 0.7126843333333334, 2.138053, 0.02862375676503981, 0.02862375676503981, 3, 1, 2.500933333333333, 7.5028, 0.129395502869887, 0.129395502869887, 1.782376005377778, 16.0413840484, 0.003703785400637656, 0.003703785400637656),
    @("sCs", "Myoc", "Fzd7", "FAPs", 3, 1, 0.7126843333333334, 2.138053, 0.02862375676503981, 0.02862375676503981, 3, 1, 7.782607, 23.347821, 0.4026634108880828, 0.4026634108880828, 5.546542081390333, 49.918878732513, 0.01152573953144177, 0.01152573953144177),
    @("sCs", "Myoc", "Fzd7", "M2", 3, 1, 0.7126843333333334, 2.138053, 0.02862375676503981, 0.02862375676503981, 3, 1, 4.565521, 13.696563, 0.2362149673420707, 0.2362149673420707, 3.253775290204334, 29.28397761183901, 0.006761359769461255, 0.006761359769461254),
    @("sCs", "Myoc", "Fzd7", "sCs", 3, 1, 0.7126843333333334, 2.138053, 0.02862375676503981, 0.02862375676503981, 3, 1, 4.478761333333333, 13.436284, 0.2317261188999595, 0.2317261188999595, 3.191943035005778, 28.727487315052, 0.006632872063499137, 0.006632872063499137),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $rowvals = $data[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowvals[$c - 1]
    }
}

Write-Host "UsedRange:" $ws.UsedRange.Address()
